# Weekly update: insert the newest "Poroto granado" (Femacal de La Calera,
# Coquimbo) price record as a new row 240, pushing the existing historical
# rows (old 240-288) down by one to 241-289.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 240 (shifts rows 240:288 -> 241:289).
$ws.Rows.Item(240).Insert()

# Populate the new row with the latest observation.
$ws.Range("A240").Value = 3
$ws.Range("B240").Value = "Femacal de La Calera"
$ws.Range("C240").Value = "Coquimbo"
$ws.Range("D240").Value = 45258
$ws.Range("E240").Value = 5
$ws.Range("F240").Value = 100112030
$ws.Range("G240").Value = "Poroto granado"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 65
$ws.Range("K240").Value = 41000
$ws.Range("L240").Value = 41000
$ws.Range("M240").Value = 41000
$ws.Range("N240").Value = "$/saco 25 kilos"
$ws.Range("O240").Value = "Provincia de Limarí"
$ws.Range("P240").Value = 1640
$ws.Range("Q240").Value = 25
$ws.Range("R240").Value = "Hortaliza"
